$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 49, pushing existing rows 49-59 down to 50-60.
$ws.Rows("49:49").Insert()

# Populate the new row 49 with the same data as the (now-shifted) row 50,
# except for the date in column D which is different (44476 instead of 44468).
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44476
$ws.Range("D49").NumberFormat = $ws.Range("D50").NumberFormat
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 100112013
$ws.Range("G49").Value = "Alcachofa"
$ws.Range("H49").Value = "Madrigal"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 300
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = 10000
$ws.Range("N49").Value = "$/caja 40 unidades"
$ws.Range("O49").Value = "Provincia del Elquí"
$ws.Range("P49").Value = 250
$ws.Range("Q49").Value = 40
$ws.Range("R49").Value = "Hortaliza"
